$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price column (D) values that changed; force text type to preserve formatting
$priceUpdates = @{
    "D2" = "67.540.22"
    "D3" = "3.243.15"
    "D5" = "577.98"
    "D6" = "179.61"
    "D8" = "0.596"
    "D9" = "3.234.85"
    "D10" = "0.130"
    "D11" = "6.77"
    "D13" = "3.797.51"
    "D15" = "28.01"
    "D16" = "67.426.19"
    "D17" = "0.0000168"
    "D18" = "3.241.56"
    "D19" = "5.81"
    "D20" = "13.38"
    "D21" = "374.23"
    "D22" = "7.59"
    "D24" = "71.21"
    "D26" = "0.0000119"
    "D27" = "9.67"
    "D31" = "5.62"
    "D32" = "22.60"
    "D34" = "1.28"
    "D35" = "6.83"
    "D36" = "164.46"
    "D37" = "1.50"
    "D38" = "0.862"
    "D40" = "6.87"
    "D41" = "26.77"
    "D42" = "363.14"
    "D43" = "2.58"
    "D45" = "2.705.57"
    "D46" = "25.81"
    "D47" = "40.43"
    "D48" = "0.0675"
    "D49" = "0.0278"
    "D51" = "0.998"
}
foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
    $cell.ClearFormats()
}

# Update Volume(1h) column (E) values
$volumeUpdates = @{
    "E2" = "  +5.36%  "
    "E3" = "  +2.90%  "
    "E4" = "  -0.08%  "
    "E5" = "  +2.99%  "
    "E6" = "  +6.80%  "
    "E7" = "  -0.06%  "
    "E8" = "  -2.28%  "
    "E9" = "  +2.58%  "
    "E10" = "  +4.77%  "
    "E11" = "  +3.56%  "
    "E12" = "  +5.01%  "
    "E13" = "  +2.73%  "
    "E14" = "  +1.16%  "
    "E15" = "  +3.43%  "
    "E16" = "  +5.22%  "
    "E17" = "  +3.13%  "
    "E18" = "  +2.69%  "
    "E19" = "  +1.67%  "
    "E20" = "  +3.88%  "
    "E21" = "  +6.72%  "
    "E22" = "  +5.91%  "
    "E23" = "  -0.33%  "
    "E24" = "  +4.72%  "
    "E25" = "  +2.18%  "
    "E26" = "  +3.87%  "
    "E27" = "  +1.24%  "
    "E28" = "  +3.17%  "
    "E29" = "  +0.40%  "
    "E30" = "  +4.99%  "
    "E31" = "  +3.38%  "
    "E32" = "  +3.52%  "
    "E33" = "  -0.04%  "
    "E34" = "  +6.72%  "
    "E35" = "  +3.94%  "
    "E36" = "  +7.00%  "
    "E37" = "  +5.13%  "
    "E38" = "  +5.86%  "
    "E39" = "  +10.06%  "
    "E40" = "  +15.19%  "
    "E41" = "  +1.83%  "
    "E42" = "  +13.45%  "
    "E43" = "  +5.48%  "
    "E44" = "  +6.18%  "
    "E45" = "  +3.22%  "
    "E46" = "  +9.13%  "
    "E47" = "  +3.00%  "
    "E48" = "  +4.32%  "
    "E49" = "  +3.21%  "
    "E50" = "  +1.17%  "
    "E51" = "  +6.76%  "
}
foreach ($addr in $volumeUpdates.Keys) {
    $ws.Range($addr).Value = $volumeUpdates[$addr]
}
